# Daily attendance processing - 2026-02-18 19:03:25 UTC
# Re-orders the comma-separated "Recorded By" values in column G so that the
# current academic session (2025/2026) is listed last, with the remaining
# entries sorted ascending before it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    "G3"  = "2022/2023, 2025/2026"
    "G22" = "2024/2025, 2025/2026"
    "G23" = "2022/2023, 2023/2024, 2025/2026"
    "G24" = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    "G27" = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    "G28" = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    "G31" = "2022/2023, 2025/2026"
    "G50" = "2024/2025, 2025/2026"
    "G51" = "2022/2023, 2023/2024, 2025/2026"
    "G52" = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    "G55" = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    "G56" = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
